$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1721.0741
$ws.Range("I40").Value = 1463.3572
$ws.Range("J40").Value = 1998.6154
$ws.Range("K40").Value = 1463.3572
$ws.Range("L40").Value = 1998.6154
$ws.Range("M40").Value = -1288.3572
$ws.Range("N40").Value = -2348.6154
$ws.Range("H64").Value = 3242.4285
$ws.Range("I64").Value = 2965.6667
$ws.Range("J64").Value = 3450
$ws.Range("K64").Value = 2965.6667
$ws.Range("L64").Value = 3450
$ws.Range("M64").Value = -2717.6667
$ws.Range("N64").Value = -3946
$ws.Range("H67").Value = 3242.4285
$ws.Range("I67").Value = 2965.6667
$ws.Range("J67").Value = 3450
$ws.Range("K67").Value = 2965.6667
$ws.Range("L67").Value = 3450
$ws.Range("M67").Value = -2107.6667
$ws.Range("N67").Value = -5166
$ws.Range("H74").Value = 4324.08
$ws.Range("I74").Value = 2800.125
$ws.Range("K74").Value = 2800.125
$ws.Range("M74").Value = -1864.125
$ws.Range("H77").Value = 4324.08
$ws.Range("I77").Value = 2800.125
$ws.Range("K77").Value = 14000.625
$ws.Range("M77").Value = -9320.625
$ws.Range("H132").Value = 32262792
$ws.Range("I132").Value = 31919524
$ws.Range("J132").Value = 33338364
$ws.Range("K132").Value = 95758572
$ws.Range("L132").Value = 100015092
$ws.Range("M132").Value = -95756042
$ws.Range("N132").Value = -100020152
$ws.Range("H138").Value = 2189.254
$ws.Range("I138").Value = 1294.7632
$ws.Range("J138").Value = 3548.88
$ws.Range("K138").Value = 3884.2896
$ws.Range("L138").Value = 10646.64
$ws.Range("M138").Value = 1255.7104
$ws.Range("N138").Value = -20926.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13557672
$ws.Range("I32").Value = 15226289
$ws.Range("K32").Value = 15226289
$ws.Range("M32").Value = -15226002
$ws.Range("H61").Value = 1216.7162
$ws.Range("I61").Value = 1272.1555
$ws.Range("J61").Value = 1130.6897
$ws.Range("K61").Value = 1272.1555
$ws.Range("L61").Value = 1130.6897
$ws.Range("M61").Value = -1060.1555
$ws.Range("N61").Value = -1554.6897
$ws.Range("H63").Value = 1677.7727
$ws.Range("I63").Value = 1663.25
$ws.Range("J63").Value = 1716.5
$ws.Range("K63").Value = 1663.25
$ws.Range("L63").Value = 1716.5
$ws.Range("M63").Value = -977.25
$ws.Range("N63").Value = -3088.5
$ws.Range("H66").Value = 1677.7727
$ws.Range("I66").Value = 1663.25
$ws.Range("J66").Value = 1716.5
$ws.Range("K66").Value = 8316.25
$ws.Range("L66").Value = 8582.5
$ws.Range("M66").Value = -4884.25
$ws.Range("N66").Value = -15446.5
$ws.Range("H110").Value = 3042
$ws.Range("I110").Value = 1215
$ws.Range("J110").Value = 4260
$ws.Range("K110").Value = 1215
$ws.Range("L110").Value = 4260
$ws.Range("M110").Value = 830
$ws.Range("N110").Value = -8350
$ws.Range("H136").Value = 1216.7162
$ws.Range("I136").Value = 1272.1555
$ws.Range("J136").Value = 1130.6897
$ws.Range("K136").Value = 3816.4665
$ws.Range("L136").Value = 3392.0691
$ws.Range("M136").Value = -1266.4665
$ws.Range("N136").Value = -8492.069100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1060007
$ws.Range("I86").Value = 3995.4167
$ws.Range("J86").Value = 2327221
$ws.Range("K86").Value = 3995.4167
$ws.Range("L86").Value = 2327221
$ws.Range("M86").Value = -2872.4167
$ws.Range("N86").Value = -2329467
$ws.Range("H89").Value = 1060007
$ws.Range("I89").Value = 3995.4167
$ws.Range("J89").Value = 2327221
$ws.Range("K89").Value = 19977.0835
$ws.Range("L89").Value = 11636105
$ws.Range("M89").Value = -14361.0835
$ws.Range("N89").Value = -11647337
$ws.Range("H105").Value = 142859550

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3213.6345
$ws.Range("I31").Value = 2230.527
$ws.Range("J31").Value = 7042.579
$ws.Range("K31").Value = 2230.527
$ws.Range("L31").Value = 7042.579
$ws.Range("M31").Value = -1935.527
$ws.Range("N31").Value = -7632.579
$ws.Range("H34").Value = 3213.6345
$ws.Range("I34").Value = 2230.527
$ws.Range("J34").Value = 7042.579
$ws.Range("K34").Value = 2230.527
$ws.Range("L34").Value = 7042.579
$ws.Range("M34").Value = -2028.527
$ws.Range("N34").Value = -7446.579
$ws.Range("H58").Value = 55556616
$ws.Range("I58").Value = 111112216
$ws.Range("J58").Value = 1019.8889
$ws.Range("K58").Value = 111112216
$ws.Range("L58").Value = 1019.8889
$ws.Range("M58").Value = -111112013
$ws.Range("N58").Value = -1425.8889
$ws.Range("H111").Value = 43351
$ws.Range("J111").Value = 43351
$ws.Range("L111").Value = 43351
$ws.Range("N111").Value = -51531
$ws.Range("H132").Value = 15157311
$ws.Range("I132").Value = 1143.1333
$ws.Range("K132").Value = 3429.3999
$ws.Range("M132").Value = -899.3998999999999
$ws.Range("H134").Value = 829.65625
$ws.Range("I134").Value = 762.3461
$ws.Range("K134").Value = 2287.0383
$ws.Range("M134").Value = 247.9616999999998
$ws.Range("H136").Value = 55556616
$ws.Range("I136").Value = 111112216
$ws.Range("J136").Value = 1019.8889
$ws.Range("K136").Value = 333336648
$ws.Range("L136").Value = 3059.6667
$ws.Range("M136").Value = -333334098
$ws.Range("N136").Value = -8159.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 806.1
$ws.Range("J131").Value = 827.06525
$ws.Range("L131").Value = 2481.19575
$ws.Range("N131").Value = -12561.19575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4166.6665
$ws.Range("I70").Value = 4107.6924
$ws.Range("K70").Value = 4107.6924
$ws.Range("M70").Value = -3837.6924
$ws.Range("H73").Value = 4166.6665
$ws.Range("I73").Value = 4107.6924
$ws.Range("K73").Value = 4107.6924
$ws.Range("M73").Value = -3171.6924
$ws.Range("H80").Value = 20005520
$ws.Range("I80").Value = 6899.75
$ws.Range("J80").Value = 100000000
$ws.Range("K80").Value = 6899.75
$ws.Range("L80").Value = 100000000
$ws.Range("M80").Value = -5901.75
$ws.Range("N80").Value = -100001996
$ws.Range("H83").Value = 20005520
$ws.Range("I83").Value = 6899.75
$ws.Range("J83").Value = 100000000
$ws.Range("K83").Value = 34498.75
$ws.Range("L83").Value = 500000000
$ws.Range("M83").Value = -29506.75
$ws.Range("N83").Value = -500009984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1348.9333
$ws.Range("I68").Value = 1273.4
$ws.Range("K68").Value = 1273.4
$ws.Range("M68").Value = -524.4000000000001
$ws.Range("H71").Value = 1348.9333
$ws.Range("I71").Value = 1273.4
$ws.Range("K71").Value = 6367
$ws.Range("M71").Value = -2623
$ws.Range("H100").Value = 3389.45
$ws.Range("I100").Value = 3123.75
$ws.Range("J100").Value = 3566.5833
$ws.Range("K100").Value = 3123.75
$ws.Range("L100").Value = 3566.5833
$ws.Range("M100").Value = -2582.75
$ws.Range("N100").Value = -4648.5833
$ws.Range("H122").Value = 3831.0754
$ws.Range("I122").Value = 4286.756
$ws.Range("J122").Value = 2274.1667
$ws.Range("K122").Value = 12860.268
$ws.Range("L122").Value = 6822.500100000001
$ws.Range("M122").Value = -10410.268
$ws.Range("N122").Value = -11722.5001
$ws.Range("H136").Value = 3001.7966
$ws.Range("I136").Value = 2725.5107
$ws.Range("K136").Value = 8176.532099999999
$ws.Range("M136").Value = -5626.532099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1750.5
$ws.Range("I122").Value = 1249.8148
$ws.Range("J122").Value = 2651.7334
$ws.Range("K122").Value = 3749.4444
$ws.Range("L122").Value = 7955.2002
$ws.Range("M122").Value = -1299.4444
$ws.Range("N122").Value = -12855.2002
$ws.Range("H136").Value = 1987.0222
$ws.Range("I136").Value = 2308.5
$ws.Range("J136").Value = 1404.3438
$ws.Range("K136").Value = 6925.5
$ws.Range("L136").Value = 4213.0314
$ws.Range("M136").Value = -4375.5
$ws.Range("N136").Value = -9313.0314
